$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2023-08-21 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-22 Tuesday", 2)

# Update the division problems in the table, addressed by row/column so that
# identical "before" and "after" values across different cells never collide.
$tbl = $d.Tables.Item(1)

function Set-CellText($row, $col, $text) {
    $cell = $tbl.Cell($row, $col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $text
}

Set-CellText 1 1 "68÷7="
Set-CellText 1 2 "15÷4="
Set-CellText 1 3 "49÷7="
Set-CellText 1 4 "69÷9="
Set-CellText 1 5 "63÷6="

Set-CellText 5 1 "20÷8="
Set-CellText 5 2 "71÷8="
Set-CellText 5 3 "24÷8="
Set-CellText 5 4 "73÷9="
Set-CellText 5 5 "23÷3="

Set-CellText 9 1 "96÷6="
Set-CellText 9 2 "28÷8="
Set-CellText 9 3 "49÷6="
Set-CellText 9 4 "92÷5="
Set-CellText 9 5 "67÷6="

Set-CellText 13 1 "69÷4="
Set-CellText 13 2 "41÷2="
Set-CellText 13 3 "70÷8="
Set-CellText 13 4 "45÷3="
Set-CellText 13 5 "13÷9="

Set-CellText 17 1 "91÷8="
Set-CellText 17 2 "40÷8="
Set-CellText 17 3 "83÷3="
Set-CellText 17 4 "93÷4="
Set-CellText 17 5 "60÷5="
